$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the OpportunityID test value in A2 (CompanyName in B2 stays "Amazon.com, Inc")
$ws.Range("A2").Value = "OPE-0003018780"
$ws.Range("B2").Value = "Amazon.com, Inc"

# Move the active selection from C11 to A3
$ws.Activate()
$ws.Range("A3").Select()
